$d = $word.ActiveDocument

# --- 1. Highlight the ROI / campaign-performance analysis questions in yellow ---
# These are the numbered "Questions:" list items that call out campaign / banner
# / placement ROI-style analysis; both the run text and the paragraph mark need
# the yellow highlight so the whole line (incl. the list bullet) is highlighted.
$highlightTargets = @(
    "What is the overall trend in user engagement throughout the campaign period?",
    "How does the size of the ad (banner) impact the number of clicks generated?",
    "Which publisher spaces (placements) yielded the highest number of displays and clicks?",
    "Is there a correlation between the cost of serving ads and the revenue generated from clicks?",
    "What is the average revenue generated per click for Company X during the campaign period?",
    "Which campaigns had the highest post-click conversion rates?",
    "Are there any specific trends or patterns in post-click sales amounts over time?",
    "How does the level of user engagement vary across different banner sizes?",
    "Which placement types result in the highest post-click conversion rates?",
    "Can we identify any seasonal patterns or fluctuations in displays and clicks throughout the campaign period?",
    "Is there a correlation between user engagement levels and the revenue generated?",
    "Are there any outliers in terms of cost, clicks, or revenue that warrant further investigation?",
    "How does the effectiveness of campaigns vary based on the size of the ad and placement type?"
)

$wdYellow = 7

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.TrimEnd("`r", "`n", [char]7)
    if ($highlightTargets -contains $text) {
        $p.Range.Font.HighlightColorIndex = $wdYellow
    }
}

# --- 2. Drop the stray _GoBack bookmark left over from the last edit point ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
